$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-10-08 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-09 Thursday", 2) | Out-Null
$d.Content.Find.Execute("850×5=", $true, $false, $false, $false, $false, $true, 1, $false, "431×4=", 2) | Out-Null
$d.Content.Find.Execute("349×9=", $true, $false, $false, $false, $false, $true, 1, $false, "390×7=", 2) | Out-Null
$d.Content.Find.Execute("946×5=", $true, $false, $false, $false, $false, $true, 1, $false, "657×6=", 2) | Out-Null
$d.Content.Find.Execute("230×9=", $true, $false, $false, $false, $false, $true, 1, $false, "747×6=", 2) | Out-Null
$d.Content.Find.Execute("811×6=", $true, $false, $false, $false, $false, $true, 1, $false, "952×2=", 2) | Out-Null
$d.Content.Find.Execute("567×6=", $true, $false, $false, $false, $false, $true, 1, $false, "444×8=", 2) | Out-Null
$d.Content.Find.Execute("446×8=", $true, $false, $false, $false, $false, $true, 1, $false, "460×8=", 2) | Out-Null
$d.Content.Find.Execute("126×7=", $true, $false, $false, $false, $false, $true, 1, $false, "480×3=", 2) | Out-Null
$d.Content.Find.Execute("564×4=", $true, $false, $false, $false, $false, $true, 1, $false, "714×9=", 2) | Out-Null
$d.Content.Find.Execute("152×4=", $true, $false, $false, $false, $false, $true, 1, $false, "604×7=", 2) | Out-Null
$d.Content.Find.Execute("198×7=", $true, $false, $false, $false, $false, $true, 1, $false, "926×5=", 2) | Out-Null
$d.Content.Find.Execute("868×6=", $true, $false, $false, $false, $false, $true, 1, $false, "341×2=", 2) | Out-Null
$d.Content.Find.Execute("780×7=", $true, $false, $false, $false, $false, $true, 1, $false, "948×6=", 2) | Out-Null
$d.Content.Find.Execute("267×2=", $true, $false, $false, $false, $false, $true, 1, $false, "421×5=", 2) | Out-Null
$d.Content.Find.Execute("817×7=", $true, $false, $false, $false, $false, $true, 1, $false, "451×7=", 2) | Out-Null
$d.Content.Find.Execute("853×9=", $true, $false, $false, $false, $false, $true, 1, $false, "970×5=", 2) | Out-Null
$d.Content.Find.Execute("829×7=", $true, $false, $false, $false, $false, $true, 1, $false, "555×3=", 2) | Out-Null
$d.Content.Find.Execute("789×7=", $true, $false, $false, $false, $false, $true, 1, $false, "429×9=", 2) | Out-Null
$d.Content.Find.Execute("398×2=", $true, $false, $false, $false, $false, $true, 1, $false, "390×7=", 2) | Out-Null
$d.Content.Find.Execute("523×3=", $true, $false, $false, $false, $false, $true, 1, $false, "593×8=", 2) | Out-Null
$d.Content.Find.Execute("339×5=", $true, $false, $false, $false, $false, $true, 1, $false, "866×9=", 2) | Out-Null
$d.Content.Find.Execute("611×4=", $true, $false, $false, $false, $false, $true, 1, $false, "688×6=", 2) | Out-Null
$d.Content.Find.Execute("886×6=", $true, $false, $false, $false, $false, $true, 1, $false, "750×5=", 2) | Out-Null
$d.Content.Find.Execute("194×6=", $true, $false, $false, $false, $false, $true, 1, $false, "645×6=", 2) | Out-Null
$d.Content.Find.Execute("453×3=", $true, $false, $false, $false, $false, $true, 1, $false, "661×7=", 2) | Out-Null
